$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two "Content" entries that actually changed text (row 10 first,
# then row 9, so the new shared-strings entries land in the same order as
# the target file: "Policy · Themes" before "Taxation · Maps").
$ws.Range("D10").Value = "Policy " + [char]0x00B7 + " Themes"
$ws.Range("D9").Value = "Taxation " + [char]0x00B7 + " Maps"

# Match the author's final active-cell selection.
$ws.Range("D10").Select()
